# Prepare bio and chemistry program list:
# - Replace the single remaining program row with the Chemical Biotechnology entry
# - Remove the other now-obsolete program rows (3-7)
# - Trim the trailing blank rows that are no longer needed (996-1000)
# - Shrink the dropdown validation range to match the remaining rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: keep "Yes" in B2, but change the program name in A2.
$ws.Range("A2").Value = "TUM_CHEMICAL_BIOTECHNOLOGY"

# Rows 3-7 no longer hold any program entries.
$ws.Range("A3:B7").Clear()

# The sheet had placeholder formatted rows all the way to 1000; trim the
# last five of those (996-1000) since they are no longer needed.
$ws.Range("A996:B1000").EntireRow.Delete()

# Keep the Yes/No dropdown validation limited to the rows that still hold data.
$ws.Range("B3:B7").Validation.Delete()

Write-Host "edit complete"
